# "Update countries & provincias Spain"
#
# The source data is a daily COVID-19 snapshot ("Pais" sheet). Between the
# "before" and "after" commits:
#   - The report timestamp in A1 was bumped from 14:52 to 15:22.
#   - Four new countries (Arabia Saudita, Emiratos Arabes Unidos, Serbia,
#     Maldivas) were (re)inserted earlier in their respective blocks of the
#     table, which pushes the countries that used to occupy those rows down
#     by one row, and every row gets that day's refreshed statistics
#     (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos
#     criticos, Muertes hoy, Muertes).
#
# Rather than physically inserting/deleting rows (which would disturb
# formatting/styles of unrelated rows), we directly write the resulting
# values for every touched cell, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" banner.
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 15:22"

# Row => final values (only columns that actually change are listed; A is
# the country name, B..H are Casos totales/Nuevos casos/Casos activos/
# Recuperados/Casos criticos/Muertes hoy/Muertes).
$rows = @(
  @{ Row = 28;  A = "Arabia Saudita";               B = 10484; C = 1122; D = 1490; E = 8891; F = 97;  G = 6;  H = 103 },
  @{ Row = 29;  A = "Chile";                        B = 10088;           D = 4338; E = 5617; F = 360;          H = 133 },
  @{ Row = 30;  A = "Ecuador";                      B = 9468;  C = 0;    D = 1061; E = 7933; F = 124; G = 0;   H = 474 },
  @{ Row = 31;  A = "Polonia";                      B = 9453;  C = 166;  D = 1133; E = 7958; F = 160; G = 2;   H = 362 },
  @{ Row = 32;  A = "Rumania";                                                     E = 6450;           G = 18; H = 469 },
  @{ Row = 37;  A = "Emiratos Arabes Unidos";        B = 7265;  C = 484;  D = 1360; E = 5862; F = 1;   G = 2;  H = 43  },
  @{ Row = 38;  A = "Noruega";                      B = 7103;  C = 25;   D = 32;   E = 6906; F = 58;  G = 0;   H = 165 },
  @{ Row = 39;  A = "Chequia";                      B = 6787;  C = 41;   D = 1311; E = 5288; F = 84;  G = 2;   H = 188 },
  @{ Row = 41;  A = "Serbia";                       B = 6630;  C = 312;  D = 870;  E = 5635; F = 108; G = 3;   H = 125 },
  @{ Row = 42;  A = "Australia";                    B = 6619;  C = 7;    D = 4258; E = 2290; F = 49;  G = 0;   H = 71  },
  @{ Row = 43;  A = "Filipinas";                    B = 6459;  C = 200;  D = 613;  E = 5418; F = 1;   G = 19;  H = 428 },
  @{ Row = 57;  A = "Argentina";                                                   D = 737;  E = 2068;         G = 2;  H = 136 },
  @{ Row = 152; A = "Maldivas";                                          C = 15;   D = 16;   E = 51;  F = 1;           H = 0   },
  @{ Row = 153; A = "San Martin (Parte Holandesa)";  B = 67;             D = 12;   E = 45;   F = 6;           H = 10  },
  @{ Row = 154; A = "Guyana";                       B = 65;             D = 9;    E = 49;   F = 4;           H = 7   },
  @{ Row = 155; A = "Islas Caimanes";                                             D = 7;    E = 53;  F = 3;           H = 1   },
  @{ Row = 156; A = "Zambia";                       B = 61;    C = 0;    D = 33;                                     H = 3   }
)

foreach ($r in $rows) {
  foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H")) {
    if ($r.ContainsKey($col)) {
      $addr = "$col$($r.Row)"
      $ws.Range($addr).Value = $r[$col]
    }
  }
}
